$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "2025-07-23 13:17:17"
$ws.Range("B8").Value = "add-user"
$ws.Range("C8").Value = "new-organization97"
$ws.Range("D8").Value = "newteam"
$ws.Range("E8").Value = "demo"
$ws.Range("F8").Value = "GokulJ17"
$ws.Range("G8").Value = "pull"

# I2 already holds the literal text "False" (not a boolean). Copy/paste its
# value into I8 so the new cell also stores "False" as text, matching how
# the rest of the column is populated, instead of Excel's usual coercion of
# a bare True/False entry into a boolean.
$ws.Range("I2").Copy()
$ws.Range("I8").PasteSpecial(-4163)
$excel.CutCopyMode = 0
